$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.977.70"
$ws.Range("E2").Value = "  +3.08%  "

# Row 3
$ws.Range("D3").Value = "3.403.88"
$ws.Range("E3").Value = "  +1.87%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "'576.82"
$ws.Range("E5").Value = "  +1.95%  "

# Row 6
$ws.Range("D6").Value = "'137.28"
$ws.Range("E6").Value = "  +5.05%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("E8").Value = "  +0.66%  "

# Row 9
$ws.Range("E9").Value = "  +0.42%  "

# Row 10
$ws.Range("E10").Value = "  +6.81%  "

# Row 11
$ws.Range("E11").Value = "  +3.97%  "

# Row 12
$ws.Range("D12").Value = "3.985.43"
$ws.Range("E12").Value = "  +1.86%  "

# Row 13
$ws.Range("E13").Value = "  +2.50%  "

# Row 15
$ws.Range("D15").Value = "3.402.34"
$ws.Range("E15").Value = "  +1.77%  "

# Row 16
$ws.Range("D16").Value = "'25.50"
$ws.Range("E16").Value = "  +3.32%  "

# Row 17
$ws.Range("D17").Value = "62.073.22"
$ws.Range("E17").Value = "  +2.97%  "

# Row 18
$ws.Range("D18").Value = "'14.24"
$ws.Range("E18").Value = "  +5.61%  "

# Row 19
$ws.Range("E19").Value = "  +3.10%  "

# Row 20
$ws.Range("E20").Value = "  +3.65%  "

# Row 21
$ws.Range("D21").Value = "'388.63"
$ws.Range("E21").Value = "  +9.37%  "

# Row 22
$ws.Range("D22").Value = "'0.569"
$ws.Range("E22").Value = "  +1.82%  "

# Row 23
$ws.Range("D23").Value = "3.546.49"
$ws.Range("E23").Value = "  +2.04%  "

# Row 24
$ws.Range("E24").Value = "  +14.27%  "

# Row 25
$ws.Range("E25").Value = "  +0.26%  "

# Row 26
$ws.Range("D26").Value = "'71.45"
$ws.Range("E26").Value = "  +2.95%  "

# Row 27
$ws.Range("D27").Value = "'7.72"
$ws.Range("E27").Value = "  +3.03%  "

# Row 28
$ws.Range("E28").Value = "  -5.83%  "

# Row 29
$ws.Range("E29").Value = "  +0.22%  "

# Row 30
$ws.Range("E30").Value = "  +4.34%  "

# Row 31
$ws.Range("D31").Value = "'0.161"
$ws.Range("E31").Value = "  +4.44%  "

# Row 32
$ws.Range("E32").Value = "  +2.21%  "

# Row 33
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.05%  "

# Row 34
$ws.Range("B34").Value = "RenzoRestakedETH"
$ws.Range("C34").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D34").Value = "3.434.57"
$ws.Range("E34").Value = "  +1.87%  "

# Row 35
$ws.Range("D35").Value = "'23.54"
$ws.Range("E35").Value = "  +2.44%  "

# Row 36
$ws.Range("D36").Value = "'5.42"
$ws.Range("E36").Value = "  +0.43%  "

# Row 37
$ws.Range("E37").Value = "  +1.11%  "

# Row 38
$ws.Range("E38").Value = "  +2.31%  "

# Row 39
$ws.Range("D39").Value = "'163.49"
$ws.Range("E39").Value = "  +2.95%  "

# Row 40
$ws.Range("D40").Value = "'0.0787"
$ws.Range("E40").Value = "  +1.85%  "

# Row 41
$ws.Range("D41").Value = "'1.79"
$ws.Range("E41").Value = "  +12.87%  "

# Row 42
$ws.Range("E42").Value = "  +4.73%  "

# Row 43
$ws.Range("E43").Value = "  +4.52%  "

# Row 44
$ws.Range("E44").Value = "  +0.03%  "

# Row 45
$ws.Range("E45").Value = "  +1.52%  "

# Row 46
$ws.Range("D46").Value = "'41.70"
$ws.Range("E46").Value = "  +2.22%  "

# Row 47
$ws.Range("D47").Value = "'24.85"
$ws.Range("E47").Value = "  +5.02%  "

# Row 48
$ws.Range("E48").Value = "  +2.34%  "

# Row 49
$ws.Range("D49").Value = "'23.35"
$ws.Range("E49").Value = "  +4.20%  "

# Row 50
$ws.Range("D50").Value = "2.371.34"
$ws.Range("E50").Value = "  +8.41%  "

# Row 51
$ws.Range("D51").Value = "'0.0264"
$ws.Range("E51").Value = "  +4.73%  "
